# Update cryptos list values (prices, volumes, coin name swap) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '47.861.69'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.509.27'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.66'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.34'
$ws.Range('E6').Value = '  +3.58%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.551'
$ws.Range('E9').Value = '  +1.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.85'
$ws.Range('E10').Value = '  +6.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0815'
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('E12').Value = '  +0.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.71'
$ws.Range('E13').Value = '  +1.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.27'
$ws.Range('E14').Value = '  +1.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.901.77'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.502.84'
$ws.Range('E16').Value = '  +0.51%  '
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '47.765.61'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('E19').Value = '  +4.24%  '
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.83'
$ws.Range('E21').Value = '  +16.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0946'
$ws.Range('E22').Value = '  +0.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.84'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '248.23'
$ws.Range('E24').Value = '  -1.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.56'
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.93'
$ws.Range('E27').Value = '  -1.21%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.140'
$ws.Range('E30').Value = '  +3.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.03'
$ws.Range('E31').Value = '  -0.67%  '
$ws.Range('E32').Value = '  +0.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.15'
$ws.Range('E33').Value = '  +2.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.39'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0788'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.69'
$ws.Range('E38').Value = '  +1.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.99'
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.112'
$ws.Range('E40').Value = '  +0.22%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.39'
$ws.Range('E41').Value = '  +5.94%  '
$ws.Range('E42').Value = '  -1.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '119.38'
$ws.Range('E43').Value = '  -1.82%  '
$ws.Range('E44').Value = '  +0.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.001.90'
$ws.Range('E45').Value = '  +1.76%  '
$ws.Range('E46').Value = '  +2.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.03'
$ws.Range('E47').Value = '  -3.39%  '
$ws.Range('E48').Value = '  +1.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.07'
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.22'
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '56.92'
$ws.Range('E51').Value = '  +3.89%  '
